$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (Chi khắc Can) text in column B with the extended wording
$ws.Range("B4").Value = "Cuộc đời của bạn nhiều điều không toại lòng. Hãy kiên định và cố gắng sẽ có quả ngọt, cẩn thận trong cả lúc thuận lợi nhất."

# Add the two new "luận đại vận" rows (order matters for shared-string table layout)
$ws.Range("A6").Value = "Âm Dương Thuận Lý"
$ws.Range("A7").Value = "Âm Dương Nghịch Lý"
$ws.Range("B7").Value = "Độ số may mắn trong cuộc đời bị giảm đi. Bạn nên kiên nhẫn chắc chắn gặt quả ngọt sau nhiều bài học bắt buộc phải có."
$ws.Range("B6").Value = "Độ số may mắn trong cuộc đời bạn được gia tăng. Đừng để ưu ái của vận mệnh khiến bạn mất ý chí cố gắng mà hãy tận dụng cơ hội để nâng cao năng lực của bản thân."

# Match the selection state recorded in the saved file
$ws.Range("U7").Select()
